# Update cryptos list - apply latest price/volume snapshot values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'26.223.77"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "'1.645.30"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D5").Value = "'216.88"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "'0.507"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("D10").Value = "'19.97"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'1.873.18"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "'1.637.69"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "'63.53"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").Value = "'26.215.56"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "'195.42"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("D21").Value = "'4.42"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'6.35"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").Value = "'15.60"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").Value = "'0.0504"
$ws.Range("E31").Value = "  +2.26%  "
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("D33").Value = "'3.24"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").Value = "'1.136.32"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "'0.553"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("D40").Value = "'0.0158"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").Value = "'100.12"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "'1.782.62"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "'56.21"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("E47").Value = "  +4.94%  "
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.73"
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.418"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  +1.39%  "
